# Commit: "Give access to camera target"
#
# The plan tracks tasks in column A/B with a Yes/No/Not Yet status in
# column C. This change marks several items (including the "Target"
# sub-item of "Camera", which grants access to the camera's target) as
# completed ("Yes").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(14, 15, 17, 18, 25, 28, 33, 39, 43, 44)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Yes"
}

# The view had scrolled so that row 31 was pinned at the top; restore the
# default scroll position (top-left back at A1).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
